## feat: add 2022-Q3 data
##
## Inserts a new "2022-Q3" sheet (fund holdings detail) right after the
## "2022-Q2" sheet tab position (i.e. right before the existing "2022-Q2"
## sheet, becoming the second tab overall, after "总计"). The previously
## second sheet ("2022-Q2") and everything after it shifts one slot to the
## right. The "总计" (summary) sheet gains a new leading data row for
## 2022-Q3 and all of its existing rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q2" detail sheet (preserves all cell
#    styles/borders/fonts) and drop the copy in right before it. The
#    duplicate becomes the brand-new "2022-Q3" sheet; the original sheet
#    (now one slot later) keeps serving as "2022-Q2".
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Overwrite the duplicated rows with the actual 2022-Q3 fund data (order
# of the two funds is swapped relative to 2022-Q2). Columns B-G are text
# in this table (e.g. "007216", "4.40") - the leading apostrophe keeps
# them as text instead of letting Excel coerce them to numbers.
$q3Sheet.Range("B2").Value = "'007216"
$q3Sheet.Range("C2").Value = "浙商中华预期高股息C"
$q3Sheet.Range("D2").Value = "'4.40"
$q3Sheet.Range("E2").Value = "'88.55"
$q3Sheet.Range("F2").Value = "'6.69"
$q3Sheet.Range("G2").Value = "'0.2944"
$q3Sheet.Range("H2").Value = 10

$q3Sheet.Range("B3").Value = "'007178"
$q3Sheet.Range("C3").Value = "浙商中华预期高股息A"
$q3Sheet.Range("D3").Value = "'2.59"
$q3Sheet.Range("E3").Value = "'88.55"
$q3Sheet.Range("F3").Value = "'6.69"
$q3Sheet.Range("G3").Value = "'0.1733"
$q3Sheet.Range("H3").Value = 10

# The leading apostrophe above is the only way to stop Excel coercing
# "007216"/"4.40"/etc. into numbers, but it also tags those cells with a
# quote-prefixed style. Column C (fund name) never needed that trick, so
# its format is still the plain, un-styled one the source table used for
# every text cell - reuse it to wipe the stray formatting back off
# columns B and D:G while leaving the freshly-typed text values in place.
$q3Sheet.Range("C2").Copy()
$q3Sheet.Range("B2").PasteSpecial(-4122)
$q3Sheet.Range("D2:G2").PasteSpecial(-4122)
$q3Sheet.Range("C3").Copy()
$q3Sheet.Range("B3").PasteSpecial(-4122)
$q3Sheet.Range("D3:G3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push every existing row down by one
#    and add a new top data row for 2022-Q3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Clone row 4's formatting into the new row 5 before touching values.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.47

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 1.05

$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q3"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.01

# ---------------------------------------------------------------------
# 3. Restore the originally-selected tab ("2021-Q3", the last sheet) -
#    copying a sheet makes the new copy active, so move the selection
#    back to where it was before this edit.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()

